$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").Value = "TJ Bamba"
$ws.Range("B62").Value = "Washington State"
$ws.Range("C62").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4596365.png&w=350&h=254"

$ws.Range("A63").Value = "Dylan Darling"
$ws.Range("B63").Value = "Washington State"
$ws.Range("C63").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105658.png&w=350&h=254"

$ws.Range("A64").Value = "Adrame Diongue"
$ws.Range("B64").Value = "Washington State"
$ws.Range("C64").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105656.png&w=350&h=254"

$ws.Range("A65").Value = "Mouhamed Gueye"
$ws.Range("B65").Value = "Washington State"
$ws.Range("C65").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4712863.png&w=350&h=254"

$ws.Range("A66").Value = "Mael Hamon-Crespin"
$ws.Range("B66").Value = "Washington State"
$ws.Range("C66").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105654.png&w=350&h=254"

$ws.Range("A67").Value = "Kymany Houinsou"
$ws.Range("B67").Value = "Washington State"
$ws.Range("C67").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105659.png&w=350&h=254"

$ws.Range("A68").Value = "Dishon Jackson"
$ws.Range("B68").Value = "Washington State"
$ws.Range("C68").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4432926.png&w=350&h=254"

$ws.Range("A69").Value = "Andrej Jakimovski"
$ws.Range("B69").Value = "Washington State"
$ws.Range("C69").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4702466.png&w=350&h=254"

$ws.Range("A70").Value = "Braden Korpela"
$ws.Range("B70").Value = "Washington State"
$ws.Range("C70").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105655.png&w=350&h=254"

$ws.Range("A71").Value = "Shae Korpela"
$ws.Range("B71").Value = "Washington State"
$ws.Range("C71").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105653.png&w=350&h=254"

$ws.Range("A72").Value = "Jabe Mullins"
$ws.Range("B72").Value = "Washington State"
$ws.Range("C72").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4701181.png&w=350&h=254"

$ws.Range("A73").Value = "Ben Olesen"
$ws.Range("B73").Value = "Washington State"
$ws.Range("C73").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4592673.png&w=350&h=254"

$ws.Range("A74").Value = "Justin Powell"
$ws.Range("B74").Value = "Washington State"
$ws.Range("C74").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4432250.png&w=350&h=254"

$ws.Range("A75").Value = "Myles Rice"
$ws.Range("B75").Value = "Washington State"
$ws.Range("C75").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4709133.png&w=350&h=254"

$ws.Range("A76").Value = "DJ Rodman"
$ws.Range("B76").Value = "Washington State"
$ws.Range("C76").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4432103.png&w=350&h=254"

$ws.Range("A77").Value = "AJ Rohosy"
$ws.Range("B77").Value = "Washington State"
$ws.Range("C77").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105657.png&w=350&h=254"

$ws.Range("A78").Value = "Carlos Rosario"
$ws.Range("B78").Value = "Washington State"
$ws.Range("C78").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4702465.png&w=350&h=254"

$ws.Range("A79").Value = "Dylan Andrews"
$ws.Range("B79").Value = "UCLA"
$ws.Range("C79").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105636.png&w=350&h=254"

$ws.Range("A80").Value = "Amari Bailey"
$ws.Range("B80").Value = "UCLA"
$ws.Range("C80").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105638.png&w=350&h=254"

$ws.Range("A81").Value = "Adem Bona"
$ws.Range("B81").Value = "UCLA"
$ws.Range("C81").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105637.png&w=350&h=254"

$ws.Range("A82").Value = "Tyger Campbell"
$ws.Range("B82").Value = "UCLA"
$ws.Range("C82").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4397128.png&w=350&h=254"

$ws.Range("A83").Value = "Abramo Canka"
$ws.Range("B83").Value = "UCLA"
$ws.Range("C83").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105635.png&w=350&h=254"

$ws.Range("A84").Value = "Jaylen Clark"
$ws.Range("B84").Value = "UCLA"
$ws.Range("C84").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4432247.png&w=350&h=254"

$ws.Range("A85").Value = "Logan Cremonesi"
$ws.Range("B85").Value = "UCLA"
$ws.Range("C85").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4702456.png&w=350&h=254"

$ws.Range("A86").Value = "Mac Etienne"
$ws.Range("B86").Value = "UCLA"
$ws.Range("C86").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4711692.png&w=350&h=254"

$ws.Range("B87").Value = "UCLA"
$ws.Range("C87").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4432848.png&w=350&h=254"

$ws.Range("A88").Value = "Evan Manjikian"
$ws.Range("B88").Value = "UCLA"
$ws.Range("C88").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105639.png&w=350&h=254"

$ws.Range("A89").Value = "Will McClendon"
$ws.Range("B89").Value = "UCLA"
$ws.Range("C89").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4433275.png&w=350&h=254"

$ws.Range("A90").Value = "Kenneth Nwuba"
$ws.Range("B90").Value = "UCLA"
$ws.Range("C90").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4397129.png&w=350&h=254"

$ws.Range("A91").Value = "Jack Seidler"
$ws.Range("B91").Value = "UCLA"
$ws.Range("C91").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/5105640.png&w=350&h=254"

$ws.Range("A92").Value = "David Singleton"
$ws.Range("B92").Value = "UCLA"
$ws.Range("C92").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4397131.png&w=350&h=254"

$ws.Range("A93").Value = "Russell Stong"
$ws.Range("B93").Value = "UCLA"
$ws.Range("C93").Value = "https://a.espncdn.com/combiner/i?img=/i/headshots/mens-college-basketball/players/full/4405598.png&w=350&h=254"

# Player name entered afterward for row 87 (UCLA)
$ws.Range("A87").Value = "Jaime Jaquez"

# Reflect the final selection/scroll position shown in the saved view
$ws.Range("A87").Select()
